$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 23450
$ws.Range("C3").Value = 705
$ws.Range("C4").Value = 10550
$ws.Range("C5").Value = 36900
$ws.Range("C6").Value = 67000
$ws.Range("C7").Value = 151500
$ws.Range("C8").Value = 46450
$ws.Range("C10").Value = 86900
$ws.Range("C11").Value = 17000
$ws.Range("C12").Value = 42450
$ws.Range("C13").Value = 9160
$ws.Range("C14").Value = 5180
$ws.Range("C15").Value = 3120
$ws.Range("C16").Value = 7750
$ws.Range("C17").Value = 72100
$ws.Range("C19").Value = 20800
$ws.Range("C20").Value = 55800
$ws.Range("C21").Value = 108500
$ws.Range("C22").Value = 75100
$ws.Range("C23").Value = 489500
$ws.Range("C24").Value = 2030
$ws.Range("C25").Value = 14700
$ws.Range("C26").Value = 30500
$ws.Range("C27").Value = 15050
$ws.Range("C28").Value = 27500
$ws.Range("C29").Value = 28850
$ws.Range("C30").Value = 35550
$ws.Range("C31").Value = 1520
$ws.Range("C32").Value = 38150
$ws.Range("C33").Value = 171500
$ws.Range("C34").Value = 3230
$ws.Range("C35").Value = 1625
$ws.Range("C36").Value = 10000
$ws.Range("C37").Value = 3405
$ws.Range("C38").Value = 3955
$ws.Range("C39").Value = 1880
$ws.Range("C40").Value = 3420
$ws.Range("C41").Value = 244000
$ws.Range("C42").Value = 952
$ws.Range("C43").Value = 87700
$ws.Range("C44").Value = 6360
$ws.Range("C45").Value = 4300
$ws.Range("C46").Value = 590
$ws.Range("C48").Value = 46150
$ws.Range("C49").Value = 18700
$ws.Range("C50").Value = 13100
$ws.Range("C51").Value = 54100

$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C52").Select()
